# "feat: adicionando planilha unificada"
# Rename the worksheet from "EXERCICIO4 - 2.3A" to "2.3A" (Excel automatically
# updates every defined name / solver_* reference that points at the sheet),
# drop the now-redundant conditional-formatting rule that only covered
# B2:H15 (the B2:H16 rule stays), and refresh the stored view state
# (top-left cell / active selection) to match what the workbook shows
# after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Rename the sheet -------------------------------------------------
# This cascades automatically to xl/workbook.xml's <sheets> entry and to
# every localSheetId="0" defined name (solver_adj, solver_lhs*, solver_opt,
# solver_rhs*, ...), which all embed the sheet name in their formula text.
$ws.Name = "2.3A"

# --- Conditional formatting --------------------------------------------
# Two rules used to cover overlapping ranges (B2:H15 and B2:H16); the
# B2:H15 one is removed, leaving only the B2:H16 rule.
$dupRange = $ws.Range("B2:H15")
for ($i = $dupRange.FormatConditions.Count; $i -ge 1; $i--) {
    $fc = $dupRange.FormatConditions.Item($i)
    if ($fc.AppliesTo.Address() -eq '$B$2:$H$15') {
        $fc.Delete()
    }
}

# --- View state ----------------------------------------------------------
# topLeftCell moves from C1 to F1 and the active selection moves from I11
# to R1.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 6   # column F
$excel.ActiveWindow.ScrollRow = 1      # row 1
$ws.Range("R1").Select()
